$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.019.63'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').Value = '2.362.28'
$ws.Range('E3').Value = '  -2.15%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '500.10'
$ws.Range('E5').Value = '  -2.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.59'
$ws.Range('E6').Value = '  -3.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.61%  '
$ws.Range('E8').Value = '  -2.60%  '
$ws.Range('D9').Value = '2.363.40'
$ws.Range('E9').Value = '  -3.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0977'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.73'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.322'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '2.781.41'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = '56.006.92'
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.35'
$ws.Range('E16').Value = '  -2.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000132'
$ws.Range('E17').Value = '  -1.62%  '
$ws.Range('D18').Value = '2.322.58'
$ws.Range('E18').Value = '  -5.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.97'
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.02'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '305.84'
$ws.Range('E21').Value = '  -2.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.25'
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.07'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').Value = '  -3.83%  '
$ws.Range('E27').Value = '  -6.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.18'
$ws.Range('E28').Value = '  -5.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.57'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = '0.0₃0709'
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('E31').Value = '  -3.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.72'
$ws.Range('E34').Value = '  -7.56%  '
$ws.Range('E35').Value = '  -5.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.58'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.17'
$ws.Range('E37').Value = '  -6.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.73'
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.04'
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.786'
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.38'
$ws.Range('E41').Value = '  -6.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '129.05'
$ws.Range('E42').Value = '  -4.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.34'
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.68'
$ws.Range('E44').Value = '  -6.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.560'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0900'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '239.30'
$ws.Range('E47').Value = '  -7.28%  '
$ws.Range('E48').Value = '  -3.00%  '
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.97'
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('E51').Value = '  -0.70%  '
